# Regenerated handoff report: the source markdown file was re-created under a
# new GUID, and new handoff xliff files (with new content hashes) were produced
# at slightly later timestamps. Update all cells/hyperlinks that embed the old
# GUID/hash/timestamps to reflect the new handoff.

$wb = $excel.ActiveWorkbook

$oldGuid = "3ecec41c-7a14-404e-8525-a610ec6ae419"
$newGuid = "d99b4360-6e36-4ed9-a9cb-3e1508614cbf"

$newHash = "562f97cea93446b0c4e764dc1ee957f425343446"

# The hyperlink target (github blob URL) itself is unchanged - only the
# visible display text is updated to the new file name.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9217715e974c67536f8c84b2099d705dc0164726/e2e/$oldGuid.md"

# --- Overview sheet: File Name (A2), Path And Name (B2, hyperlinked), ---
# --- Latest HO Xliff Generate Date (G2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
# Replace the hyperlink on B2 so its cached display text matches the new
# path, while keeping the same underlying target URL.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkUrl, "", "", "e2e\$newGuid.md")
$wsOverview.Range("G2").Value = "2016-10-24 09:43:01"

# --- zh-cn sheet: Source File Name (A2, hyperlinked), Latest Handoff File ---
# --- (G2), Latest Handoff Datetime (H2) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkUrl, "", "", "$newGuid.md")
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-10-24 09:42:49"

# --- de-de sheet: Source File Name (A2, hyperlinked), Latest Handoff File ---
# --- (G2). Latest Handoff Datetime (H2) is untouched by this handoff. ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkUrl, "", "", "$newGuid.md")
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
